# Apply new header ordering / block-order data layout to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header row (A1:F1)
$headers = @("kitchens_1", "bedrooms_1", "living_rooms_1", "living_rooms_2", "kitchens_2", "bedrooms_2")
for ($col = 1; $col -le 6; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# New one-hot block-order data for rows 2-7 (columns A-F)
$data = @(
    @(0, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1, 0),
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$i][$col - 1]
    }
}
